$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 237.5
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(6, 8).Value = 1067.9166
$ws.Cells.Item(6, 9).Value = 301.5
$ws.Cells.Item(6, 11).Value = 904.5
$ws.Cells.Item(6, 13).Value = -792.5
$ws.Cells.Item(17, 8).Value = 837.2679000000001
$ws.Cells.Item(17, 10).Value = 842.37036
$ws.Cells.Item(17, 12).Value = 2527.11108
$ws.Cells.Item(17, 14).Value = -2863.11108
$ws.Cells.Item(19, 8).Value = 3531.5
$ws.Cells.Item(19, 9).Value = 999
$ws.Cells.Item(19, 11).Value = 999
$ws.Cells.Item(19, 13).Value = -824
$ws.Cells.Item(32, 8).Value = 9815.625
$ws.Cells.Item(32, 9).Value = 9110.299999999999
$ws.Cells.Item(32, 11).Value = 9110.299999999999
$ws.Cells.Item(32, 13).Value = -8784.299999999999
$ws.Cells.Item(39, 8).Value = 512.25
$ws.Cells.Item(39, 9).Value = 500
$ws.Cells.Item(39, 10).Value = 516.3333
$ws.Cells.Item(39, 11).Value = 1500
$ws.Cells.Item(39, 12).Value = 1548.9999
$ws.Cells.Item(39, 13).Value = -1204
$ws.Cells.Item(39, 14).Value = -2140.9999
$ws.Cells.Item(40, 8).Value = 4888.222
$ws.Cells.Item(40, 9).Value = 4499.1665
$ws.Cells.Item(40, 10).Value = 5666.3335
$ws.Cells.Item(40, 11).Value = 4499.1665
$ws.Cells.Item(40, 12).Value = 5666.3335
$ws.Cells.Item(40, 13).Value = -4324.1665
$ws.Cells.Item(40, 14).Value = -6016.3335
$ws.Cells.Item(41, 8).Value = 1547.2727
$ws.Cells.Item(41, 9).Value = 447.33334
$ws.Cells.Item(41, 10).Value = 1959.75
$ws.Cells.Item(41, 11).Value = 447.33334
$ws.Cells.Item(41, 12).Value = 1959.75
$ws.Cells.Item(41, 13).Value = -7.333340000000021
$ws.Cells.Item(41, 14).Value = -2839.75
$ws.Cells.Item(62, 8).Value = 4030.795
$ws.Cells.Item(62, 9).Value = 3582.5588
$ws.Cells.Item(62, 11).Value = 3582.5588
$ws.Cells.Item(62, 13).Value = -2958.5588
$ws.Cells.Item(65, 8).Value = 4030.795
$ws.Cells.Item(65, 9).Value = 3582.5588
$ws.Cells.Item(65, 11).Value = 17912.794
$ws.Cells.Item(65, 13).Value = -14792.794
$ws.Cells.Item(68, 8).Value = 59999
$ws.Cells.Item(68, 10).Value = 59999
$ws.Cells.Item(68, 12).Value = 59999
$ws.Cells.Item(68, 14).Value = -61497
$ws.Cells.Item(71, 8).Value = 59999
$ws.Cells.Item(71, 10).Value = 59999
$ws.Cells.Item(71, 12).Value = 179997
$ws.Cells.Item(71, 14).Value = -187485
$ws.Cells.Item(107, 8).Value = 1600.75
$ws.Cells.Item(107, 9).Value = 1115.2142
$ws.Cells.Item(107, 11).Value = 1115.2142
$ws.Cells.Item(107, 13).Value = 804.7858000000001
$ws.Cells.Item(137, 8).Value = 15311.435
$ws.Cells.Item(137, 9).Value = 18837.166
$ws.Cells.Item(137, 11).Value = 56511.49800000001
$ws.Cells.Item(137, 13).Value = -53961.49800000001
$ws.Cells.Item(141, 8).Value = 1849.5
$ws.Cells.Item(141, 9).Value = 1849.5
$ws.Cells.Item(141, 11).Value = 5548.5
$ws.Cells.Item(141, 13).Value = -368.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 12792.75
$ws.Cells.Item(5, 9).Value = 14565.286
$ws.Cells.Item(5, 11).Value = 14565.286
$ws.Cells.Item(5, 13).Value = -14453.286
$ws.Cells.Item(32, 8).Value = 23807.791
$ws.Cells.Item(32, 9).Value = 25735.318
$ws.Cells.Item(32, 11).Value = 25735.318
$ws.Cells.Item(32, 13).Value = -25448.318
$ws.Cells.Item(61, 8).Value = 5616.9644
$ws.Cells.Item(61, 9).Value = 1159.7778
$ws.Cells.Item(61, 11).Value = 1159.7778
$ws.Cells.Item(61, 13).Value = -947.7778000000001
$ws.Cells.Item(97, 8).Value = 1760.8636
$ws.Cells.Item(97, 9).Value = 1462
$ws.Cells.Item(97, 10).Value = 2119.5
$ws.Cells.Item(97, 11).Value = 1462
$ws.Cells.Item(97, 12).Value = 2119.5
$ws.Cells.Item(97, 13).Value = -966
$ws.Cells.Item(97, 14).Value = -3111.5
$ws.Cells.Item(110, 8).Value = 1499
$ws.Cells.Item(110, 9).Value = 1332.6666
$ws.Cells.Item(110, 10).Value = 1998
$ws.Cells.Item(110, 11).Value = 1332.6666
$ws.Cells.Item(110, 12).Value = 1998
$ws.Cells.Item(110, 13).Value = 712.3334
$ws.Cells.Item(110, 14).Value = -6088
$ws.Cells.Item(132, 8).Value = 1571.0862
$ws.Cells.Item(132, 9).Value = 985.2353000000001
$ws.Cells.Item(132, 11).Value = 2955.7059
$ws.Cells.Item(132, 13).Value = -425.7058999999999
$ws.Cells.Item(136, 8).Value = 5616.9644
$ws.Cells.Item(136, 9).Value = 1159.7778
$ws.Cells.Item(136, 11).Value = 3479.3334
$ws.Cells.Item(136, 13).Value = -929.3334000000004
$ws.Cells.Item(139, 8).Value = 132999.8
$ws.Cells.Item(139, 10).Value = 132999.8
$ws.Cells.Item(139, 12).Value = 132999.8
$ws.Cells.Item(139, 14).Value = -143279.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 12792.75
$ws.Cells.Item(4, 9).Value = 14565.286
$ws.Cells.Item(4, 11).Value = 14565.286
$ws.Cells.Item(4, 13).Value = -14450.286
$ws.Cells.Item(64, 8).Value = 3332.6667
$ws.Cells.Item(64, 9).Value = 3797
$ws.Cells.Item(64, 10).Value = 3100.5
$ws.Cells.Item(64, 11).Value = 3797
$ws.Cells.Item(64, 12).Value = 3100.5
$ws.Cells.Item(64, 13).Value = -3572
$ws.Cells.Item(64, 14).Value = -3550.5
$ws.Cells.Item(67, 8).Value = 3332.6667
$ws.Cells.Item(67, 9).Value = 3797
$ws.Cells.Item(67, 10).Value = 3100.5
$ws.Cells.Item(67, 11).Value = 3797
$ws.Cells.Item(67, 12).Value = 3100.5
$ws.Cells.Item(67, 13).Value = -3017
$ws.Cells.Item(67, 14).Value = -4660.5
$ws.Cells.Item(134, 8).Value = 2280.1875
$ws.Cells.Item(134, 9).Value = 1519.2609
$ws.Cells.Item(134, 10).Value = 4224.778
$ws.Cells.Item(134, 11).Value = 4557.7827
$ws.Cells.Item(134, 12).Value = 12674.334
$ws.Cells.Item(134, 13).Value = -2022.7827
$ws.Cells.Item(134, 14).Value = -17744.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5265417.5
$ws.Cells.Item(31, 9).Value = 7143937.5
$ws.Cells.Item(31, 11).Value = 7143937.5
$ws.Cells.Item(31, 13).Value = -7143642.5
$ws.Cells.Item(34, 8).Value = 5265417.5
$ws.Cells.Item(34, 9).Value = 7143937.5
$ws.Cells.Item(34, 11).Value = 7143937.5
$ws.Cells.Item(34, 13).Value = -7143735.5
$ws.Cells.Item(62, 8).Value = 4468.533
$ws.Cells.Item(62, 9).Value = 4703.1665
$ws.Cells.Item(62, 10).Value = 3530
$ws.Cells.Item(62, 11).Value = 4703.1665
$ws.Cells.Item(62, 12).Value = 3530
$ws.Cells.Item(62, 13).Value = -4079.1665
$ws.Cells.Item(62, 14).Value = -4778
$ws.Cells.Item(65, 8).Value = 4468.533
$ws.Cells.Item(65, 9).Value = 4703.1665
$ws.Cells.Item(65, 10).Value = 3530
$ws.Cells.Item(65, 11).Value = 23515.8325
$ws.Cells.Item(65, 12).Value = 17650
$ws.Cells.Item(65, 13).Value = -20395.8325
$ws.Cells.Item(65, 14).Value = -23890
$ws.Cells.Item(105, 8).Value = 25051.438
$ws.Cells.Item(105, 9).Value = 32811.5
$ws.Cells.Item(105, 11).Value = 32811.5
$ws.Cells.Item(105, 13).Value = -31064.5
$ws.Cells.Item(122, 8).Value = 1874.6316
$ws.Cells.Item(122, 9).Value = 1945.7858
$ws.Cells.Item(122, 10).Value = 1675.4
$ws.Cells.Item(122, 11).Value = 5837.357400000001
$ws.Cells.Item(122, 12).Value = 5026.200000000001
$ws.Cells.Item(122, 13).Value = -3387.357400000001
$ws.Cells.Item(122, 14).Value = -9926.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 10500000
$ws.Cells.Item(11, 9).Value = 10500000
$ws.Cells.Item(11, 11).Value = 31500000
$ws.Cells.Item(11, 13).Value = -31499860
$ws.Cells.Item(17, 8).Value = 294.2
$ws.Cells.Item(17, 9).Value = 290.33334
$ws.Cells.Item(17, 11).Value = 871.0000200000001
$ws.Cells.Item(17, 13).Value = -702.0000200000001
$ws.Cells.Item(23, 8).Value = 362.33334
$ws.Cells.Item(23, 10).Value = 477.77777
$ws.Cells.Item(23, 12).Value = 1433.33331
$ws.Cells.Item(23, 14).Value = -1903.33331
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).ClearContents()
$ws.Cells.Item(75, 14).ClearContents()
$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 13).ClearContents()
$ws.Cells.Item(78, 14).ClearContents()
$ws.Cells.Item(121, 8).Value = 64032.117
$ws.Cells.Item(121, 10).Value = 1718.8334
$ws.Cells.Item(121, 12).Value = 5156.5002
$ws.Cells.Item(121, 14).Value = -7776.5002
$ws.Cells.Item(136, 8).Value = 1848.1428
$ws.Cells.Item(136, 9).Value = 1848.1428
$ws.Cells.Item(136, 11).Value = 5544.428400000001
$ws.Cells.Item(136, 13).Value = -444.4284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 69490
$ws.Cells.Item(102, 9).Value = 77552
$ws.Cells.Item(102, 10).Value = 4994
$ws.Cells.Item(102, 11).Value = 77552
$ws.Cells.Item(102, 12).Value = 4994
$ws.Cells.Item(102, 13).Value = -75930
$ws.Cells.Item(102, 14).Value = -8238
$ws.Cells.Item(126, 8).Value = 3495.7693
$ws.Cells.Item(126, 9).Value = 2799.375
$ws.Cells.Item(126, 10).Value = 4610
$ws.Cells.Item(126, 11).Value = 8398.125
$ws.Cells.Item(126, 12).Value = 13830
$ws.Cells.Item(126, 13).Value = -5928.125
$ws.Cells.Item(126, 14).Value = -18770
$ws.Cells.Item(135, 8).Value = 100000
$ws.Cells.Item(135, 10).Value = 100000
$ws.Cells.Item(135, 12).Value = 100000
$ws.Cells.Item(135, 14).Value = -110140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2325.7334
$ws.Cells.Item(7, 9).Value = 2274
$ws.Cells.Item(7, 10).Value = 3050
$ws.Cells.Item(7, 11).Value = 2274
$ws.Cells.Item(7, 12).Value = 3050
$ws.Cells.Item(7, 13).Value = -2162
$ws.Cells.Item(7, 14).Value = -3274
$ws.Cells.Item(46, 8).Value = 3845.353
$ws.Cells.Item(46, 10).Value = 5966.6
$ws.Cells.Item(46, 12).Value = 5966.6
$ws.Cells.Item(46, 14).Value = -6342.6
$ws.Cells.Item(126, 8).Value = 2325.7334
$ws.Cells.Item(126, 9).Value = 2274
$ws.Cells.Item(126, 10).Value = 3050
$ws.Cells.Item(126, 11).Value = 6822
$ws.Cells.Item(126, 12).Value = 9150
$ws.Cells.Item(126, 13).Value = -4352
$ws.Cells.Item(126, 14).Value = -14090
$ws.Cells.Item(132, 8).Value = 1670.875
$ws.Cells.Item(132, 9).Value = 757.1429000000001
$ws.Cells.Item(132, 10).Value = 2950.1
$ws.Cells.Item(132, 11).Value = 2271.4287
$ws.Cells.Item(132, 12).Value = 8850.299999999999
$ws.Cells.Item(132, 13).Value = 258.5712999999996
$ws.Cells.Item(132, 14).Value = -13910.3
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 119000
$ws.Cells.Item(140, 10).Value = 119000
$ws.Cells.Item(140, 12).Value = 119000
$ws.Cells.Item(140, 14).Value = -129360

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2492.0715
$ws.Cells.Item(132, 9).Value = 1939.1
$ws.Cells.Item(132, 10).Value = 3874.5
$ws.Cells.Item(132, 11).Value = 5817.299999999999
$ws.Cells.Item(132, 12).Value = 11623.5
$ws.Cells.Item(132, 13).Value = -3287.299999999999
$ws.Cells.Item(132, 14).Value = -16683.5
